# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.932.08"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "1.622.13"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("E4").Value = "  +0.53%  "
$ws.Range("D5").Value = "'213.50"
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("D6").Value = "'0.501"
$ws.Range("E6").Value = "  -1.14%  "
$ws.Range("E7").Value = "  +0.53%  "
$ws.Range("E8").Value = "  -2.54%  "
$ws.Range("D9").Value = "'0.0617"
$ws.Range("E9").Value = "  -3.40%  "
$ws.Range("D10").Value = "'18.10"
$ws.Range("E10").Value = "  -7.57%  "
$ws.Range("D11").Value = "'0.0789"
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").Value = "1.848.22"
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("D13").Value = "1.621.97"
$ws.Range("E13").Value = "  -2.14%  "
$ws.Range("E14").Value = "  -2.33%  "
$ws.Range("D15").Value = "'0.522"
$ws.Range("E15").Value = "  -3.96%  "
$ws.Range("D16").Value = "25.922.44"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "'61.08"
$ws.Range("E17").Value = "  -3.52%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.0₃0733"
$ws.Range("E18").Value = "  -3.90%  "
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").Value = "'190.47"
$ws.Range("E20").Value = "  -2.18%  "
$ws.Range("D21").Value = "'4.22"
$ws.Range("E21").Value = "  -2.97%  "
$ws.Range("D22").Value = "'9.53"
$ws.Range("E22").Value = "  -3.83%  "
$ws.Range("D23").Value = "'6.04"
$ws.Range("E23").Value = "  -2.50%  "
$ws.Range("D25").Value = "'143.21"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").Value = "'1.76"
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("D28").Value = "'6.69"
$ws.Range("E28").Value = "  -2.73%  "
$ws.Range("D29").Value = "'15.15"
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("D31").Value = "'0.0479"
$ws.Range("E31").Value = "  -3.54%  "
$ws.Range("E32").Value = "  -4.73%  "
$ws.Range("D33").Value = "'3.10"
$ws.Range("E33").Value = "  -5.89%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'2.40"
$ws.Range("E34").Value = "  -2.54%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'1.48"
$ws.Range("E35").Value = "  -3.20%  "
$ws.Range("D36").Value = "1.126.40"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").Value = "'0.843"
$ws.Range("E37").Value = "  -6.70%  "
$ws.Range("E39").Value = "  -4.91%  "
$ws.Range("E40").Value = "  -2.21%  "
$ws.Range("D41").Value = "'97.75"
$ws.Range("E41").Value = "  -1.33%  "
$ws.Range("D42").Value = "'0.769"
$ws.Range("E42").Value = "  -3.23%  "
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("D44").Value = "'5.17"
$ws.Range("E44").Value = "  -5.56%  "
$ws.Range("E45").Value = "  -3.37%  "
$ws.Range("D46").Value = "'54.35"
$ws.Range("E46").Value = "  -3.98%  "
$ws.Range("D47").Value = "'0.0522"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").Value = "'7.46"
$ws.Range("E51").Value = "  -3.72%  "
